$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-09-05 Friday" "2025-09-06 Saturday"

Replace-Text "619÷7=" "633÷7="
Replace-Text "291÷6=" "250÷6="
Replace-Text "957÷6=" "407÷4="
Replace-Text "171÷7=" "523÷4="
Replace-Text "405÷9=" "885÷2="
Replace-Text "776÷8=" "762÷5="
Replace-Text "875÷4=" "972÷4="
Replace-Text "161÷6=" "528÷6="
Replace-Text "541÷3=" "493÷7="
Replace-Text "724÷2=" "134÷8="
Replace-Text "489÷2=" "475÷8="
Replace-Text "584÷5=" "625÷9="
Replace-Text "983÷6=" "458÷6="
Replace-Text "452÷2=" "984÷9="
Replace-Text "618÷9=" "568÷5="
Replace-Text "142÷6=" "345÷2="
Replace-Text "708÷5=" "387÷9="
Replace-Text "808÷2=" "900÷7="
Replace-Text "379÷4=" "858÷4="
Replace-Text "582÷5=" "512÷6="
Replace-Text "456÷4=" "268÷2="
Replace-Text "307÷7=" "296÷2="
Replace-Text "854÷5=" "623÷3="
Replace-Text "900÷4=" "389÷3="
Replace-Text "183÷5=" "713÷6="

Write-Output "Done replacing text"
